$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update "Förändrad" column C for rows 2-16: 46073 -> 46074
for ($r = 2; $r -le 16; $r++) {
    $ws.Cells.Item($r, 3).Value = 46074
}

# Rows 8-15 got re-sorted; apply the new Beteckning/Datum/Area values per row
$newData = @{
    8  = @("A 45370-2022", 44844.6397337963, 2.7)
    9  = @("A 23677-2023", 45077, 0.6)
    10 = @("A 2253-2022", 44578, 0.3)
    11 = @("A 58926-2025", 45986, 3.1)
    12 = @("A 23678-2023", 45077, 1.4)
    13 = @("A 50277-2024", 45600.60440972223, 0.5)
    14 = @("A 50530-2024", 45601.56424768519, 0.7)
    15 = @("A 50538-2024", 45601.57153935185, 0.8)
}

foreach ($r in $newData.Keys) {
    $vals = $newData[$r]
    $ws.Cells.Item($r, 1).Value = $vals[0]
    $ws.Cells.Item($r, 2).Value = $vals[1]
    $ws.Cells.Item($r, 7).Value = $vals[2]
}
